# Rename the "paymentplannames" sheet to "paymentplan" and restructure its
# contents: the payment-plan lookup rows now carry plan-economics columns
# (C:I) on the first data row, the lookup names shift from "M-kopa test
# 10..14" to "M-kopa test 26..30", and the stale extra rows (7..16) are
# removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paymentplannames")
$ws.Name = "paymentplan"

# --- Clear out the old extra rows (7:16) so nothing lingers below row 6 ---
$ws.Rows("7:16").Delete() | Out-Null

# --- Row 1: headers ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Display Name (Marketing):"
$ws.Range("C1").Value = "cashprice"
$ws.Range("D1").Value = "loanDepoit"
$ws.Range("E1").Value = "Freeusage"
$ws.Range("F1").Value = "costpercredit"
$ws.Range("G1").Value = "days"
$ws.Range("H1").Value = "loantotalprice"
$ws.Range("I1").Value = "Approvernotes"

# --- Row 2 (note: "Test" is introduced into the shared-string table ahead
#     of the renamed lookup values below, matching the author's edit order) ---
$ws.Range("I2").Value = "Test"
$ws.Range("A2").Value = "M-kopa test 26"
$ws.Range("B2").Value = "M-kopa test 26"
$ws.Range("C2").Value = 63000
$ws.Range("D2").Value = 2999
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 55
$ws.Range("G2").Value = 365
$ws.Range("H2").Value = 65000

# --- Row 3 ---
$ws.Range("A3").Value = "M-kopa test 27"
$ws.Range("B3").Value = "M-kopa test 27"

# --- Row 4 ---
$ws.Range("A4").Value = "M-kopa test 28"
$ws.Range("B4").Value = "M-kopa test 28"

# --- Row 5 ---
$ws.Range("A5").Value = "M-kopa test 29"
$ws.Range("B5").Value = "M-kopa test 29"

# --- Row 6 ---
$ws.Range("A6").Value = "M-kopa test 30"
$ws.Range("B6").Value = "M-kopa test 30"

# Match the author's saved selection on the sheet (A2:B2, active cell A2)
$ws.Activate()
$ws.Range("A2:B2").Select()
